$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 8 ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Cells.Item(8, 1).Value = "Retour status"
$ws.Cells.Item(8, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(8, 4).Value = "Retour / Terugbetaling"
$ws.Cells.Item(8, 6).Value = "2025-08-28 18:15:13"
$ws.Cells.Item(8, 7).Value = "Nee"
$ws.Cells.Item(8, 8).Value = "Ja"
$ws.Cells.Item(8, 9).Value = "Nee"
$ws.Cells.Item(8, 10).Value = "Nee"

# Extend the conditional formatting ranges from row 7 to row 8
$ranges = @("D2:D7", "G2:G7", "H2:H7", "I2:I7", "J2:J7")
foreach ($r in $ranges) {
  $col = $r.Substring(0, 1)
  $fcs = $ws.Range($r).FormatConditions
  $fc = $fcs.Item(1)
  $newRange = $ws.Range("$($col)2:$($col)8")
  $fc.ModifyAppliesToRange($newRange)
}

# --- Dashboard sheet: bump the count in B2 from 6 to 7 ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 7
